$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 22 values (D22, F22)
$ws.Range("D22").Value = 9.990842018928157
$ws.Range("F22").Value = 19.99084201892816

# Add new rows 26-33 for 2025-02-07 and 2025-02-08
$data = @(
    @("2025-02-07", "abs_activity", 10, 10, 0, 20),
    @("2025-02-07", "rel_activity", 10, 5.135171668018383, 0, 15.13517166801838),
    @("2025-02-07", "abs_sleep", 6.866666666666667, 8.266666666666667, 0, 15.13333333333333),
    @("2025-02-07", "rel_sleep", 0, 0, 0, 0),
    @("2025-02-08", "abs_activity", 0, 0, 0, 0),
    @("2025-02-08", "rel_activity", 10, 0, 0, 10),
    @("2025-02-08", "abs_sleep", 0, 0, 0, 0),
    @("2025-02-08", "rel_sleep", 0, 0, 0, 0)
)

$row = 26
foreach ($entry in $data) {
    # Column A holds date-like text (e.g. "2025-02-07") that must stay a
    # literal text value rather than being auto-converted to a date serial.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 1).ClearFormats()
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row++
}
